$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("201×6=", $true, $false, $false, $false, $false, $true, 1, $false, "371×3=", 2)
$null = $d.Content.Find.Execute("257×6=", $true, $false, $false, $false, $false, $true, 1, $false, "107×4=", 2)
$null = $d.Content.Find.Execute("411×7=", $true, $false, $false, $false, $false, $true, 1, $false, "751×8=", 2)
$null = $d.Content.Find.Execute("191×2=", $true, $false, $false, $false, $false, $true, 1, $false, "514×5=", 2)
$null = $d.Content.Find.Execute("580×2=", $true, $false, $false, $false, $false, $true, 1, $false, "561×8=", 2)
$null = $d.Content.Find.Execute("643×6=", $true, $false, $false, $false, $false, $true, 1, $false, "500×2=", 2)
$null = $d.Content.Find.Execute("333×9=", $true, $false, $false, $false, $false, $true, 1, $false, "995×2=", 2)
$null = $d.Content.Find.Execute("771×5=", $true, $false, $false, $false, $false, $true, 1, $false, "778×6=", 2)
$null = $d.Content.Find.Execute("466×9=", $true, $false, $false, $false, $false, $true, 1, $false, "810×7=", 2)
$null = $d.Content.Find.Execute("551×3=", $true, $false, $false, $false, $false, $true, 1, $false, "963×6=", 2)
$null = $d.Content.Find.Execute("503×4=", $true, $false, $false, $false, $false, $true, 1, $false, "261×9=", 2)
$null = $d.Content.Find.Execute("282×6=", $true, $false, $false, $false, $false, $true, 1, $false, "958×4=", 2)
$null = $d.Content.Find.Execute("863×7=", $true, $false, $false, $false, $false, $true, 1, $false, "127×9=", 2)
$null = $d.Content.Find.Execute("371×5=", $true, $false, $false, $false, $false, $true, 1, $false, "773×3=", 2)
$null = $d.Content.Find.Execute("684×3=", $true, $false, $false, $false, $false, $true, 1, $false, "747×3=", 2)
$null = $d.Content.Find.Execute("690×5=", $true, $false, $false, $false, $false, $true, 1, $false, "351×9=", 2)
$null = $d.Content.Find.Execute("841×6=", $true, $false, $false, $false, $false, $true, 1, $false, "139×8=", 2)
$null = $d.Content.Find.Execute("482×9=", $true, $false, $false, $false, $false, $true, 1, $false, "520×9=", 2)
$null = $d.Content.Find.Execute("170×2=", $true, $false, $false, $false, $false, $true, 1, $false, "894×6=", 2)
$null = $d.Content.Find.Execute("192×6=", $true, $false, $false, $false, $false, $true, 1, $false, "873×7=", 2)
$null = $d.Content.Find.Execute("524×2=", $true, $false, $false, $false, $false, $true, 1, $false, "566×3=", 2)
$null = $d.Content.Find.Execute("119×9=", $true, $false, $false, $false, $false, $true, 1, $false, "536×2=", 2)
$null = $d.Content.Find.Execute("507×4=", $true, $false, $false, $false, $false, $true, 1, $false, "465×7=", 2)
$null = $d.Content.Find.Execute("316×4=", $true, $false, $false, $false, $false, $true, 1, $false, "773×8=", 2)
$null = $d.Content.Find.Execute("316×9=", $true, $false, $false, $false, $false, $true, 1, $false, "376×7=", 2)
